$p = $ppt.ActivePresentation

# Slide 1 ("Forst litt lek med tall...") had its title placeholder
# removed entirely.
$s = $p.Slides.Item(1)

# Locate the title placeholder ("Tittel 1" / ctrTitle) by name and
# remember its placeholder type so we can find the regenerated stub
# afterwards regardless of what PowerPoint happens to name it.
$title = $s.Shapes.Item("Tittel 1")
$titleType = $title.PlaceholderFormat.Type

# Deleting a layout placeholder first clears it back to an empty,
# freshly (re)created placeholder stub - this mirrors PowerPoint's own
# delete-placeholder behaviour, and matches what the authors'
# change-tracking log for this edit recorded: the original shape is
# marked deleted and a new stub shape is added-then-deleted within the
# very same action. Deleting it a second time removes that stub
# entirely so no title shape remains on the slide at all.
$title.Delete()

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq $titleType) {
        $shp.Delete()
    }
}
